# Add two new columns (I: "I0", J: "IF") to the right of the existing
# data (which runs through column H), matching the header style used by
# the other header cells in row 1 (copy H1's formatting, then overwrite
# the value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers - copy H1's style (bold, bordered, centered) onto I1/J1, then set text
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# Data rows 2-33: column I and J values
$data = @{
    2  = @(8, 9)
    3  = @(1, 3)
    4  = @(1, 7)
    5  = @(1, 5)
    6  = @(1, 5)
    7  = @(1, 5)
    8  = @(1, 4)
    9  = @(1, 5)
    10 = @(1, 5)
    11 = @(1, 5)
    12 = @(1, 5)
    13 = @(1, 4)
    14 = @(1, 6)
    15 = @(1, 6)
    16 = @(1, 6)
    17 = @(1, 7)
    18 = @(1, 6)
    19 = @(1, 6)
    20 = @(1, 6)
    21 = @(1, 6)
    22 = @(1, 6)
    23 = @(1, 5)
    24 = @(1, 5)
    25 = @(1, 6)
    26 = @(1, 6)
    27 = @(1, 7)
    28 = @(1, 6)
    29 = @(1, 6)
    30 = @(1, 4)
    31 = @(1, 3)
    32 = @(1, 2)
    33 = @(1, 1)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
